$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 13:14, shifting the old rows 13-21 down to 15-23.
# Formatting/height travels with the shifted rows automatically.
$ws.Rows("13:14").Insert(-4121)

# Fix "Objetivos:" body text (row 10)
$ws.Range("B10").Value = 'Introduzir os conceitos fundamentais da ciência administração e de configurações de uma organização.'
$ws.Range("C10").Value = 'Introduzir os conceitos fundamentais da ciência administração e de configurações de uma organização.'

# New row 13: first professor under "Docentes responsaveis:"
# Copy column formatting from row 15 (A/B/C = styles 1/2/3) so the new cells
# land on the same cellXfs as the rest of the sheet instead of a fresh style.
$ws.Range("A15:C15").Copy()
$ws.Range("A13:C13").PasteSpecial(-4122)
$ws.Range("A13").Clear()
$ws.Range("B13").Value = '11079086 - Herlandí de Souza Andrade'
$ws.Range("C13").Value = '11079086 - Herlandí de Souza Andrade'

# New row 14: second professor
$ws.Range("A15:C15").Copy()
$ws.Range("A14:C14").PasteSpecial(-4122)
$ws.Range("A14").Clear()
$ws.Range("B14").Value = '5840560 - Marco Antonio Carvalho Pereira'
$ws.Range("C14").Value = '5840560 - Marco Antonio Carvalho Pereira'

# Row 15 (was 13): "Programa resumido:" body text
$ws.Range("B15").Value = "1. Conceitos Fundamentais de Administração.`n2. Noções Básicas de Estratégia."
$ws.Range("C15").Value = "1. Conceitos Fundamentais de Administração.`n2. Noções Básicas de Estratégia."

# Row 17 (was 15): "Programa:" body text
$ws.Range("B17").Value = "1. Teoria Geral de Administração: Histórico. Visão de Taylor. Escola clássica de administração.`n2. Conceitos básicos de Estratégia, Análise SWOT, Balanced ScoreCard, Mapas Estratégicos"
$ws.Range("C17").Value = "1. Teoria Geral de Administração: Histórico. Visão de Taylor. Escola clássica de administração.`n2. Conceitos básicos de Estratégia, Análise SWOT, Balanced ScoreCard, Mapas Estratégicos"

# Row 20 (was 18): "Metodo:" body text
$ws.Range("B20").Value = 'Aulas expositivas e dialogadas; dinâmicas, projetos e trabalhos em grupo; exercícios individuais; e, seminários, debates e palestras.'
$ws.Range("C20").Value = 'Aulas expositivas e dialogadas; dinâmicas, projetos e trabalhos em grupo; exercícios individuais; e, seminários, debates e palestras.'

# Row 21 (was 19): "Criterio:" body text
$ws.Range("B21").Value = 'Média Aritmética dos Projetos, Trabalhos e Exercícios realizados no decorrer da disciplina, considerando as questões relativas às Competências (Conhecimento, Habilidade e Atitude) desenvolvidas.'
$ws.Range("C21").Value = 'Média Aritmética dos Projetos, Trabalhos e Exercícios realizados no decorrer da disciplina, considerando as questões relativas às Competências (Conhecimento, Habilidade e Atitude) desenvolvidas.'

# Row 22 (was 20): "Norma de recuperacao:" body text
$ws.Range("B22").Value = 'NF = (MF + PR)/2, onde MF é a média final da avaliação e PR é uma prova de recuperação.'
$ws.Range("C22").Value = 'NF = (MF + PR)/2, onde MF é a média final da avaliação e PR é uma prova de recuperação.'

# Row 23 (was 21): "Bibliografia:" body text
$ws.Range("B23").Value = "CHIAVENATO, I; SAPIRO, A. Planejamento Estratégico. Rio de Janeiro. Campus, 2004 `nCOLLINS, J.C.; PORRAS, J. I. Feitas para Durar: Práticas bem-sucedidas de empresas visionárias. 9ª Ed.  Rio de Janeiro. Rocco, 2007 `nHERRERO, E. Balanced Scorecard e a Gestão Estratégica. Rio de Janeiro. Campus, 2005. `nKAPLAN, R; NORTON, D. Kaplan e Norton na Prática. Rio de Janeiro. Campus, 2004 `nKAPLAN, R; NORTON, D. A Estratégia em Ação: Balanced Scorecard. Rio de Janeiro. Campus, 1997 `nKAPLAN, R; NORTON, D. Mapas Estratégicos. Rio de Janeiro. Campus, 2004 `nTZU, S. A Arte da Guerra (Edição Completa). São Paulo. WMF Martins Fontes, 2009."
$ws.Range("C23").Value = "CHIAVENATO, I; SAPIRO, A. Planejamento Estratégico. Rio de Janeiro. Campus, 2004 `nCOLLINS, J.C.; PORRAS, J. I. Feitas para Durar: Práticas bem-sucedidas de empresas visionárias. 9ª Ed.  Rio de Janeiro. Rocco, 2007 `nHERRERO, E. Balanced Scorecard e a Gestão Estratégica. Rio de Janeiro. Campus, 2005. `nKAPLAN, R; NORTON, D. Kaplan e Norton na Prática. Rio de Janeiro. Campus, 2004 `nKAPLAN, R; NORTON, D. A Estratégia em Ação: Balanced Scorecard. Rio de Janeiro. Campus, 1997 `nKAPLAN, R; NORTON, D. Mapas Estratégicos. Rio de Janeiro. Campus, 2004 `nTZU, S. A Arte da Guerra (Edição Completa). São Paulo. WMF Martins Fontes, 2009."

# Column A now only spans column 1 (was redundantly min=1 max=2); narrowing it
# this way leaves column A width untouched and column B keeps its own width/style.
$ws.Range("B:B").ColumnWidth = $ws.Range("B:B").ColumnWidth
